# Update workbook to match target: 19 odds corrections in rows 2-4,
# plus one new match row (row 7) appended to the Sheet1 table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Correct existing odds values in rows 2-4 ---
$ws.Range("L2").Value = 2.25
$ws.Range("O2").Value = 1.44
$ws.Range("P2").Value = 2.63
$ws.Range("Q2").Value = 2.35
$ws.Range("R2").Value = 1.57
$ws.Range("U2").Value = 2.38
$ws.Range("V2").Value = 1.53
$ws.Range("Z2").Value = 81
$ws.Range("AB2").Value = 67
$ws.Range("AY2").Value = 26
$ws.Range("BD2").Value = 126

$ws.Range("O3").Value = 1.53
$ws.Range("P3").Value = 2.38
$ws.Range("Q3").Value = 2.7
$ws.Range("R3").Value = 1.44

$ws.Range("O4").Value = 1.18
$ws.Range("P4").Value = 4.5
$ws.Range("Q4").Value = 1.62
$ws.Range("R4").Value = 2.25

# --- Append new row 7 (new match: Dep. La Coruna vs Eibar) ---
# Force B7 to stay as literal text (it looks like a date) ...
$ws.Range("B7").NumberFormat = "@"
$ws.Range("A7").Value = "Wj3sKyQp"
$ws.Range("B7").Value = "11/11/2024"
$ws.Range("C7").Value = "16:30"
$ws.Range("D7").Value = "SPAIN - LALIGA2"
$ws.Range("E7").Value = "Dep. La Coruna"
$ws.Range("F7").Value = "Eibar"
$ws.Range("G7").Value = 1.75
$ws.Range("H7").Value = 3.5
$ws.Range("I7").Value = 4.75
$ws.Range("J7").Value = 2.4
$ws.Range("K7").Value = 2.1
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 10
$ws.Range("O7").Value = 1.3
$ws.Range("P7").Value = 3.4
$ws.Range("Q7").Value = 2.03
$ws.Range("R7").Value = 1.83
$ws.Range("S7").Value = 1.44
$ws.Range("T7").Value = 2.63
$ws.Range("U7").Value = 1.83
$ws.Range("V7").Value = 1.83
$ws.Range("W7").Value = 6.5
$ws.Range("X7").Value = 8
$ws.Range("Y7").Value = 8.5
$ws.Range("Z7").Value = 13
$ws.Range("AA7").Value = 15
$ws.Range("AB7").Value = 29
$ws.Range("AC7").Value = 9
$ws.Range("AD7").Value = 6.5
$ws.Range("AE7").Value = 17
$ws.Range("AF7").Value = 51
$ws.Range("AG7").Value = 301
$ws.Range("AH7").Value = 12
$ws.Range("AI7").Value = 23
$ws.Range("AJ7").Value = 15
$ws.Range("AK7").Value = 51
$ws.Range("AL7").Value = 41
$ws.Range("AM7").Value = 41
$ws.Range("AN7").Value = 3.75
$ws.Range("AO7").Value = 9.5
$ws.Range("AP7").Value = 21
$ws.Range("AQ7").Value = 34
$ws.Range("AR7").Value = 51
$ws.Range("AS7").Value = 151
$ws.Range("AT7").Value = 2.63
$ws.Range("AU7").Value = 8.5
$ws.Range("AV7").Value = 51
$ws.Range("AW7").Value = 6.5
$ws.Range("AX7").Value = 26
$ws.Range("AY7").Value = 34
$ws.Range("AZ7").Value = 81
$ws.Range("BA7").Value = 126
$ws.Range("BB7").Value = 251
$ws.Range("BC7").Value = 81
$ws.Range("BD7").Value = 81

# ... then drop the temporary number-format override so the cell
# keeps the default (unstyled) look of the rest of the data rows.
$ws.Range("B7").ClearFormats()

"done"